# "added harvester and experiment design"
#
# Sheet1 header row: A=harvestDate, B=harvester, C=biosampleNumber,
# D=experimentDesign, E=experimentObservations, F=strain, G=genotype, ...
#
# For every data row (2-13) fill in the newly-tracked metadata:
#   - harvester (B)         -> "S.GISH"      (replaces placeholder "Retrofitted_480")
#   - experimentDesign (D)  -> "90minuteInduction" (new column, was blank)
#   - strain (F)            -> "KN99alpha"   (new column, was blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
    $ws.Cells.Item($r, 6).Value = "KN99alpha"
}
